$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1009.6667
$ws.Range("I9").Value = 1299.5555
$ws.Range("J9").Value = 574.8333
$ws.Range("K9").Value = 1299.5555
$ws.Range("L9").Value = 574.8333
$ws.Range("M9").Value = -1130.5555
$ws.Range("N9").Value = -912.8333
$ws.Range("H33").Value = 16667803
$ws.Range("H76").Value = 4182.1875
$ws.Range("J76").Value = 4896.5
$ws.Range("L76").Value = 4896.5
$ws.Range("N76").Value = -5526.5
$ws.Range("H79").Value = 4182.1875
$ws.Range("J79").Value = 4896.5
$ws.Range("L79").Value = 4896.5
$ws.Range("N79").Value = -7080.5
$ws.Range("H86").Value = 10901.516
$ws.Range("I86").Value = 8550.666999999999
$ws.Range("K86").Value = 8550.666999999999
$ws.Range("M86").Value = -7427.666999999999
$ws.Range("H89").Value = 10901.516
$ws.Range("I89").Value = 8550.666999999999
$ws.Range("K89").Value = 42753.335
$ws.Range("M89").Value = -37137.335
$ws.Range("H92").Value = 1239.5883
$ws.Range("I92").Value = 1728.5
$ws.Range("K92").Value = 1728.5
$ws.Range("M92").Value = -480.5
$ws.Range("H94").Value = 1338.4445
$ws.Range("I94").Value = 1338.4445
$ws.Range("K94").Value = 1338.4445
$ws.Range("M94").Value = -887.4445000000001
$ws.Range("H103").Value = 1083.8235
$ws.Range("I103").Value = 736.5
$ws.Range("J103").Value = 1392.5555
$ws.Range("K103").Value = 2209.5
$ws.Range("L103").Value = 4177.666499999999
$ws.Range("M103").Value = -1623.5
$ws.Range("N103").Value = -5349.666499999999
$ws.Range("H112").Value = 2334.6333
$ws.Range("J112").Value = 2474.6072
$ws.Range("L112").Value = 7423.821599999999
$ws.Range("N112").Value = -9639.821599999999
$ws.Range("H113").Value = 15303.9
$ws.Range("J113").Value = 14120.8
$ws.Range("L113").Value = 14120.8
$ws.Range("N113").Value = -20628.8
$ws.Range("H116").Value = 4805.95
$ws.Range("I116").Value = 4832
$ws.Range("J116").Value = 4701.75
$ws.Range("K116").Value = 4832
$ws.Range("L116").Value = 4701.75
$ws.Range("M116").Value = -1390
$ws.Range("N116").Value = -11585.75
$ws.Range("H123").Value = 105463.164
$ws.Range("J123").Value = 105463.164
$ws.Range("L123").Value = 105463.164
$ws.Range("N123").Value = -115263.164
$ws.Range("H127").Value = 334330.34
$ws.Range("I127").Value = 455269.2
$ws.Range("K127").Value = 1365807.6
$ws.Range("M127").Value = -1360847.6
$ws.Range("H132").Value = 16274.857
$ws.Range("I132").Value = 15720
$ws.Range("K132").Value = 47160
$ws.Range("M132").Value = -44630
$ws.Range("H138").Value = 3987.6667
$ws.Range("J138").Value = 2845.625
$ws.Range("L138").Value = 8536.875
$ws.Range("N138").Value = -18816.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8929.714
$ws.Range("I32").Value = 3073.4932
$ws.Range("K32").Value = 3073.4932
$ws.Range("M32").Value = -2786.4932
$ws.Range("H61").Value = 9359.045
$ws.Range("I61").Value = 4022.6428
$ws.Range("J61").Value = 18148.412
$ws.Range("K61").Value = 4022.6428
$ws.Range("L61").Value = 18148.412
$ws.Range("M61").Value = -3810.6428
$ws.Range("N61").Value = -18572.412
$ws.Range("H74").Value = 13334.527
$ws.Range("I74").Value = 1781.2609
$ws.Range("K74").Value = 1781.2609
$ws.Range("M74").Value = -907.2609
$ws.Range("H76").Value = 36162.668
$ws.Range("J76").Value = 36162.668
$ws.Range("L76").Value = 36162.668
$ws.Range("N76").Value = -36838.668
$ws.Range("H77").Value = 13334.527
$ws.Range("I77").Value = 1781.2609
$ws.Range("K77").Value = 8906.3045
$ws.Range("M77").Value = -4538.3045
$ws.Range("H79").Value = 36162.668
$ws.Range("J79").Value = 36162.668
$ws.Range("L79").Value = 36162.668
$ws.Range("N79").Value = -38502.668
$ws.Range("H102").Value = 23559.4
$ws.Range("I102").Value = 3955
$ws.Range("J102").Value = 199999
$ws.Range("K102").Value = 3955
$ws.Range("L102").Value = 199999
$ws.Range("M102").Value = -2333
$ws.Range("N102").Value = -203243
$ws.Range("H110").Value = 4066.0476
$ws.Range("I110").Value = 4781.9375
$ws.Range("J110").Value = 1775.2
$ws.Range("K110").Value = 4781.9375
$ws.Range("L110").Value = 1775.2
$ws.Range("M110").Value = -2736.9375
$ws.Range("N110").Value = -5865.2
$ws.Range("H130").Value = 77923.875
$ws.Range("J130").Value = 77923.875
$ws.Range("L130").Value = 77923.875
$ws.Range("N130").Value = -87963.875
$ws.Range("H136").Value = 9359.045
$ws.Range("I136").Value = 4022.6428
$ws.Range("J136").Value = 18148.412
$ws.Range("K136").Value = 12067.9284
$ws.Range("L136").Value = 54445.236
$ws.Range("M136").Value = -9517.928400000001
$ws.Range("N136").Value = -59545.236

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 96000
$ws.Range("J59").Value = 96000
$ws.Range("L59").Value = 96000
$ws.Range("N59").Value = -97694
$ws.Range("H80").Value = 2375.625
$ws.Range("I80").Value = 700
$ws.Range("J80").Value = 2615
$ws.Range("K80").Value = 700
$ws.Range("L80").Value = 2615
$ws.Range("M80").Value = 298
$ws.Range("N80").Value = -4611
$ws.Range("H83").Value = 2375.625
$ws.Range("I83").Value = 700
$ws.Range("J83").Value = 2615
$ws.Range("K83").Value = 3500
$ws.Range("L83").Value = 13075
$ws.Range("M83").Value = 1492
$ws.Range("N83").Value = -23059
$ws.Range("H99").Value = 10113.167
$ws.Range("I99").Value = 665.5714
$ws.Range("K99").Value = 665.5714
$ws.Range("M99").Value = 832.4286
$ws.Range("H107").Value = 2428.32
$ws.Range("I107").Value = 2090.2632
$ws.Range("K107").Value = 2090.2632
$ws.Range("M107").Value = -170.2631999999999
$ws.Range("H110").Value = 40000.668
$ws.Range("J110").Value = 40000.668
$ws.Range("L110").Value = 40000.668
$ws.Range("N110").Value = -48180.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2035.5714
$ws.Range("I16").Value = 1671.4286
$ws.Range("K16").Value = 1671.4286
$ws.Range("M16").Value = -1384.4286
$ws.Range("H31").Value = 19631.473
$ws.Range("I31").Value = 9003.143
$ws.Range("K31").Value = 9003.143
$ws.Range("M31").Value = -8708.143
$ws.Range("H34").Value = 19631.473
$ws.Range("I34").Value = 9003.143
$ws.Range("K34").Value = 9003.143
$ws.Range("M34").Value = -8801.143
$ws.Range("H44").Value = 25000
$ws.Range("I44").Value = 25000
$ws.Range("K44").Value = 25000
$ws.Range("M44").Value = -24558
$ws.Range("H99").Value = 6711.846
$ws.Range("I99").Value = 3520.4443
$ws.Range("K99").Value = 3520.4443
$ws.Range("M99").Value = -2022.4443
$ws.Range("H105").Value = 13090
$ws.Range("I105").Value = 20494
$ws.Range("J105").Value = 750
$ws.Range("K105").Value = 20494
$ws.Range("L105").Value = 750
$ws.Range("M105").Value = -18747
$ws.Range("N105").Value = -4244
$ws.Range("H107").Value = 1526.0769
$ws.Range("I107").Value = 1055.2142
$ws.Range("K107").Value = 1055.2142
$ws.Range("M107").Value = 864.7858000000001
$ws.Range("H113").Value = 2035.5714
$ws.Range("I113").Value = 1671.4286
$ws.Range("K113").Value = 1671.4286
$ws.Range("M113").Value = 498.5714
$ws.Range("H126").Value = 6711.846
$ws.Range("I126").Value = 3520.4443
$ws.Range("K126").Value = 10561.3329
$ws.Range("M126").Value = -8091.332900000001
$ws.Range("H132").Value = 14335.637
$ws.Range("I132").Value = 4150.2
$ws.Range("K132").Value = 12450.6
$ws.Range("M132").Value = -9920.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 621.2308
$ws.Range("I13").Value = 95.833336
$ws.Range("J13").Value = 1071.5714
$ws.Range("K13").Value = 287.500008
$ws.Range("L13").Value = 3214.7142
$ws.Range("M13").Value = -119.500008
$ws.Range("N13").Value = -3550.7142
$ws.Range("H47").Value = 940
$ws.Range("I47").Value = 1750
$ws.Range("K47").Value = 5250
$ws.Range("M47").Value = -4819
$ws.Range("H69").Value = 13500
$ws.Range("J69").Value = 17500
$ws.Range("L69").Value = 52500
$ws.Range("N69").Value = -54122
$ws.Range("H72").Value = 13500
$ws.Range("J72").Value = 17500
$ws.Range("L72").Value = 157500
$ws.Range("N72").Value = -165612
$ws.Range("H80").Value = 12545.363
$ws.Range("J80").Value = 19199.8
$ws.Range("L80").Value = 57599.39999999999
$ws.Range("N80").Value = -59471.39999999999
$ws.Range("H83").Value = 12545.363
$ws.Range("J83").Value = 19199.8
$ws.Range("L83").Value = 172798.2
$ws.Range("N83").Value = -182158.2
$ws.Range("H86").Value = 824.5833
$ws.Range("I86").Value = 794.44446
$ws.Range("K86").Value = 2383.33338
$ws.Range("M86").Value = -1197.33338
$ws.Range("H87").Value = 15000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H89").Value = 824.5833
$ws.Range("I89").Value = 794.44446
$ws.Range("K89").Value = 7150.00014
$ws.Range("M89").Value = -1222.00014
$ws.Range("H90").Value = 15000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H97").Value = 10067.857
$ws.Range("I97").Value = 491.66666
$ws.Range("J97").Value = 17250
$ws.Range("K97").Value = 1474.99998
$ws.Range("L97").Value = 51750
$ws.Range("M97").Value = -978.9999800000001
$ws.Range("N97").Value = -52742
$ws.Range("H107").Value = 3473700
$ws.Range("J107").Value = 6251580
$ws.Range("L107").Value = 18754740
$ws.Range("N107").Value = -18758580
$ws.Range("H113").Value = 1611.0834
$ws.Range("I113").Value = 1325.5
$ws.Range("J113").Value = 1753.875
$ws.Range("K113").Value = 3976.5
$ws.Range("L113").Value = 5261.625
$ws.Range("M113").Value = -1806.5
$ws.Range("N113").Value = -9601.625
$ws.Range("H131").Value = 1499.6
$ws.Range("J131").Value = 1499.6
$ws.Range("L131").Value = 4498.799999999999
$ws.Range("N131").Value = -14578.8
$ws.Range("H133").Value = 6832
$ws.Range("J133").Value = 9749.5
$ws.Range("L133").Value = 29248.5
$ws.Range("N133").Value = -39368.5
$ws.Range("H139").Value = 11980.263
$ws.Range("I139").Value = 17362.5
$ws.Range("J139").Value = 6000
$ws.Range("K139").Value = 52087.5
$ws.Range("L139").Value = 18000
$ws.Range("M139").Value = -46947.5
$ws.Range("N139").Value = -28280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H80").Value = 6724.788
$ws.Range("I80").Value = 3914.7827
$ws.Range("J80").Value = 13187.8
$ws.Range("K80").Value = 3914.7827
$ws.Range("L80").Value = 13187.8
$ws.Range("M80").Value = -2916.7827
$ws.Range("N80").Value = -15183.8
$ws.Range("H83").Value = 6724.788
$ws.Range("I83").Value = 3914.7827
$ws.Range("J83").Value = 13187.8
$ws.Range("K83").Value = 19573.9135
$ws.Range("L83").Value = 65939
$ws.Range("M83").Value = -14581.9135
$ws.Range("N83").Value = -75923
$ws.Range("H102").Value = 3737.077
$ws.Range("I102").Value = 3793.5217
$ws.Range("J102").Value = 3304.3333
$ws.Range("K102").Value = 3793.5217
$ws.Range("L102").Value = 3304.3333
$ws.Range("M102").Value = -2171.5217
$ws.Range("N102").Value = -6548.3333
$ws.Range("H113").Value = 46701.1
$ws.Range("I113").Value = 61754.266
$ws.Range("K113").Value = 61754.266
$ws.Range("M113").Value = -59584.266
$ws.Range("H132").Value = 7298.794
$ws.Range("I132").Value = 4294.5806
$ws.Range("K132").Value = 12883.7418
$ws.Range("M132").Value = -10353.7418

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3046.4614
$ws.Range("J16").Value = 1868.4
$ws.Range("L16").Value = 1868.4
$ws.Range("N16").Value = -2208.4
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H68").Value = 2103964.5
$ws.Range("I68").Value = 1696.8
$ws.Range("J68").Value = 2979909.2
$ws.Range("K68").Value = 1696.8
$ws.Range("L68").Value = 2979909.2
$ws.Range("M68").Value = -947.8
$ws.Range("N68").Value = -2981407.2
$ws.Range("H69").Value = 1000000000
$ws.Range("J69").Value = 1000000000
$ws.Range("L69").Value = 1000000000
$ws.Range("N69").Value = -1000001622
$ws.Range("H71").Value = 2103964.5
$ws.Range("I71").Value = 1696.8
$ws.Range("J71").Value = 2979909.2
$ws.Range("K71").Value = 8484
$ws.Range("L71").Value = 14899546
$ws.Range("M71").Value = -4740
$ws.Range("N71").Value = -14907034
$ws.Range("H72").Value = 1000000000
$ws.Range("J72").Value = 1000000000
$ws.Range("L72").Value = 3000000000
$ws.Range("N72").Value = -3000008112
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H82").Value = 5448.7
$ws.Range("I82").Value = 5031.6665
$ws.Range("J82").Value = 6699.8
$ws.Range("K82").Value = 5031.6665
$ws.Range("L82").Value = 6699.8
$ws.Range("M82").Value = -4670.6665
$ws.Range("N82").Value = -7421.8
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H85").Value = 5448.7
$ws.Range("I85").Value = 5031.6665
$ws.Range("J85").Value = 6699.8
$ws.Range("K85").Value = 5031.6665
$ws.Range("L85").Value = 6699.8
$ws.Range("M85").Value = -3783.6665
$ws.Range("N85").Value = -9195.799999999999
$ws.Range("H100").Value = 2792.4546
$ws.Range("I100").Value = 2041.9231
$ws.Range("J100").Value = 3876.5557
$ws.Range("K100").Value = 2041.9231
$ws.Range("L100").Value = 3876.5557
$ws.Range("M100").Value = -1500.9231
$ws.Range("N100").Value = -4958.5557
$ws.Range("H132").Value = 1298437
$ws.Range("I132").Value = 4797.4546
$ws.Range("J132").Value = 4460667
$ws.Range("K132").Value = 14392.3638
$ws.Range("L132").Value = 13382001
$ws.Range("M132").Value = -11862.3638
$ws.Range("N132").Value = -13387061

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4662.8887
$ws.Range("J14").Value = 4852.2856
$ws.Range("L14").Value = 4852.2856
$ws.Range("N14").Value = -5188.2856
$ws.Range("H41").Value = 20940.5
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H107").Value = 1119.2858
$ws.Range("I107").Value = 1093.0588
$ws.Range("J107").Value = 1159.8182
$ws.Range("K107").Value = 3279.1764
$ws.Range("L107").Value = 3479.4546
$ws.Range("M107").Value = -1359.1764
$ws.Range("N107").Value = -7319.4546
$ws.Range("H126").Value = 18516.65
$ws.Range("I126").Value = 21549.117
$ws.Range("J126").Value = 1332.6666
$ws.Range("K126").Value = 64647.351
$ws.Range("L126").Value = 3997.9998
$ws.Range("M126").Value = -62177.351
$ws.Range("N126").Value = -8937.9998
$ws.Range("H140").Value = 160587.67
$ws.Range("J140").Value = 160587.67
$ws.Range("L140").Value = 160587.67
$ws.Range("N140").Value = -170947.67
